$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new configuration row (Name / Value / Description) to the
# config table on Sheet1.
$ws.Range("A75").Value = "CustomerUrlRegex"
$ws.Range("B75").Value = ".*person\/(\d*)\/permissions"
$ws.Range("C75").Value = "Regex to extract the customer number from the RP Url"

# Grow the Excel Table (Table1) so the new row becomes part of it.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:C75"))

# Update the view/selection state to match where the user left off.
$ws.Range("C77").Select()
